$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "done" status for the TwitterDAL column (E) on rows that
# previously had no value there.
$ws.Range("E6").Value = "done"
$ws.Range("E8").Value = "done"
$ws.Range("E10").Value = "done"

# Move the active selection to E13, matching the final state of the file.
$ws.Range("E13").Select()
